# Update the division problems in the table to the new set of values.
# Using table cell addressing (row/column) avoids ambiguity from duplicate
# text (e.g. "692÷8=" appears twice in the original document but must map
# to two different new values depending on its position).

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; Old="881÷7="; New="319÷9="},
    @{Row=1;  Col=2; Old="352÷2="; New="659÷5="},
    @{Row=1;  Col=3; Old="692÷8="; New="832÷9="},
    @{Row=1;  Col=4; Old="939÷8="; New="816÷5="},
    @{Row=1;  Col=5; Old="758÷5="; New="350÷2="},

    @{Row=5;  Col=1; Old="983÷9="; New="936÷9="},
    @{Row=5;  Col=2; Old="511÷7="; New="216÷6="},
    @{Row=5;  Col=3; Old="692÷8="; New="666÷9="},
    @{Row=5;  Col=4; Old="388÷3="; New="426÷3="},
    @{Row=5;  Col=5; Old="364÷9="; New="982÷5="},

    @{Row=9;  Col=1; Old="586÷7="; New="817÷7="},
    @{Row=9;  Col=2; Old="690÷7="; New="504÷6="},
    @{Row=9;  Col=3; Old="598÷5="; New="109÷2="},
    @{Row=9;  Col=4; Old="912÷8="; New="200÷4="},
    @{Row=9;  Col=5; Old="389÷5="; New="656÷7="},

    @{Row=13; Col=1; Old="716÷3="; New="513÷7="},
    @{Row=13; Col=2; Old="293÷8="; New="905÷6="},
    @{Row=13; Col=3; Old="828÷9="; New="829÷2="},
    @{Row=13; Col=4; Old="846÷7="; New="830÷5="},
    @{Row=13; Col=5; Old="480÷9="; New="258÷4="},

    @{Row=17; Col=1; Old="608÷8="; New="727÷4="},
    @{Row=17; Col=2; Old="269÷8="; New="806÷3="},
    @{Row=17; Col=3; Old="576÷5="; New="155÷4="},
    @{Row=17; Col=4; Old="895÷5="; New="875÷6="},
    @{Row=17; Col=5; Old="298÷5="; New="796÷3="}
)

foreach ($item in $replacements) {
    $cell = $table.Cell($item.Row, $item.Col)
    $range = $cell.Range
    # Use wdReplaceOne (1) rather than wdReplaceAll (2) so the replacement
    # stays confined to this cell's range, even when the same search text
    # (e.g. "692÷8=") occurs elsewhere in the table.
    $range.Find.Execute($item.Old, $true, $false, $false, $false, $false, `
                         $true, 0, $false, $item.New, 1)
}

Write-Host "Updated $($replacements.Count) cells."
